# The sheet's two "sib1" estimate/se/n column triples (J:L = sib1_sex0_*,
# M:O = sib1_sex1_*) were reordered so that sib1_sex1_* now comes first
# (J:L) and sib1_sex0_* moves to M:O. This is the same swap for every row,
# including row 1 (the header row, whose shared-string labels follow the
# same J<->M, K<->N, L<->O pattern).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 1
$lastRow = 93

for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($offset = 0; $offset -lt 3; $offset++) {
        $colLeft = 10 + $offset   # J, K, L
        $colRight = 13 + $offset  # M, N, O

        $leftCell = $ws.Cells.Item($r, $colLeft)
        $rightCell = $ws.Cells.Item($r, $colRight)

        $leftValue = $leftCell.Value()
        $rightValue = $rightCell.Value()

        $leftCell.Value = $rightValue
        $rightCell.Value = $leftValue
    }
}

# The active selection moved from F58 to C59:C90.
$ws.Range("C59:C90").Select()
